$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.598.42'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.918.33'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.41%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4874'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2912'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06747'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '111.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.31%  '
$ws.Range('D12').Value = '1.926.63'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07583'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.376'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6756'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '293.99'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').Value = '30.582.56'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('E18').Value = '  +1.19%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007579'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.12%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.539'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('D22').Value = '2.165.38'
$ws.Range('E22').Value = '  -0.09%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.457'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.30%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.508'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.46'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.111'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1074'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.437'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +3.08%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.158'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.14%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.091'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05027'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7412'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.141'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02036'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.702'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.686'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.022'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '109.99'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.05%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4469'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8656'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.33%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.887'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '70.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.65%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.001'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.281'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.43%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.310'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.62%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '48.30'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1234'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.2537'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.45%  '
